# PlayerPerformance_4487.xlsx — add "ODI Bowling Extra" sheet, trim empty
# cells from "ODI Batting Extra" rows that have no batting-position data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "ODI Batting Extra" — clear B:E for the rows that never had a
#    batting position recorded (the player didn't bat in that match).
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$emptyRows = @(2, 4, 10, 14, 17, 19)
foreach ($r in $emptyRows) {
    $rowRange = $battingExtra.Range($battingExtra.Cells.Item($r, 2), $battingExtra.Cells.Item($r, 5))
    $rowRange.ClearContents()
}

# ---------------------------------------------------------------------
# 2) Add a new "ODI Bowling Extra" sheet after "ODI Batting Extra",
#    mirroring its layout, and fill it with the scraped data.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Copy the header formatting (bold, border, centered) from the sibling sheet.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$data = New-Object 'object[,]' 20,3
$rows = @(
    @("3906", "0", ""),
    @("3910", "", ""),
    @("3911", "0", "30.00%"),
    @("4011", "0", ""),
    @("4038", "0", ""),
    @("4044", "0", ""),
    @("4054", "", ""),
    @("4055", "0", "10.00%"),
    @("4058", "", ""),
    @("4059", "0", ""),
    @("4060", "", ""),
    @("4099", "0", ""),
    @("4124", "1", ""),
    @("4231", "0", ""),
    @("4232", "0", ""),
    @("4375", "0", ""),
    @("4449", "0", ""),
    @("4450", "", ""),
    @("4451", "0", "10.00%"),
    @("4463", "0", "10.00%")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $data[$i, 0] = $rows[$i][0]
    $data[$i, 1] = $rows[$i][1]
    $data[$i, 2] = $rows[$i][2]
}

$range = $bowlingExtra.Range("A2:C21")
$range.NumberFormat = "@"
$range.Value = $data
